$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused trailing columns (X:AQ) that used to hold the
# repeated 1Pair-A..MaxUnique headers for the extra HKL blocks.
$ws.Range("X1:AQ1").EntireColumn.Delete()

# Rename the existing "HexGrid" rows (16-19) to the new "Holden" scheme.
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("B17").Value = "Holden5"
$ws.Range("B18").Value = "Holden10"
$ws.Range("B19").Value = "Holden15"

# Append the "HexGrid" rows as new rows 20-23, keeping the same 21 data
# columns (C:W) of 1s.
$hexNames = @("HexGrid-90degTilt2.5degRes", "HexGrid-90degTilt5degRes", "HexGrid-90degTilt10degRes", "HexGrid-90degTilt15degRes")
for ($i = 0; $i -lt 4; $i++) {
    $r = 20 + $i
    $ws.Cells.Item($r, 1).Value = 18 + $i
    $ws.Cells.Item($r, 2).Value = $hexNames[$i]
    for ($c = 3; $c -le 23; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
    # Column A keeps the bold/bordered "index" style used by the rest of
    # the table (same as column A in every other row).
    $ws.Range("A19").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
